# Work on the DUA sheet: add a "Field"/"Content" header row above the
# existing Restrictions/Terms table, shifting everything down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DUA")
$srcHeader = $wb.Worksheets.Item("participants_info").Cells.Item(1, 1)

$ws.Rows.Item(1).Insert() | Out-Null
$ws.Cells.Item(1, 1).Value = "Field"
$ws.Cells.Item(1, 2).Value = "Content"

# Reuse the existing bold/red heading format (style already present in the
# workbook) instead of letting Excel synthesize a brand-new style entry.
$srcHeader.Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

# The drop-down validation on the "Restrictions" cell (now B2) pointed at
# $D$1:$D$5; after the insert it needs to reference $D$2:$D$6.
$rng = $ws.Range("B2")
$rng.Validation.Modify(3, 1, 1, '$D$2:$D$6') | Out-Null

# Make the DUA sheet the active tab/selected cell, matching the author's
# saved view state.
$ws.Activate() | Out-Null
$ws.Range("A5").Select() | Out-Null

# Unrelated fix noted in the commit message: the dataset_info sheet should
# no longer be the tab shown as selected (handled automatically above,
# since activating DUA clears tabSelected on the previously active sheet).
